$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Header row relabel ---
$ws1.Range("B2").Value = "part_name"
$ws1.Range("C2").Value = "no_work_order"
$ws1.Range("D2").Value = "customer"
$ws1.Range("E2").Value = "prod_date"
$ws1.Range("F2").Value = "quantity_perbox"
$ws1.Range("G2").Value = "total_order"
$ws1.Range("H2").Value = "total_box"
$ws1.Range("I2").Value = "supplier"

# --- Data row update ---
$ws1.Range("B3").Value = "GARNISH RR BUMPER LWR (IPR)"

$c3 = $ws1.Range("C3")
$c3.Font.Name = "Arial"
$c3.Font.Size = 10

$ws1.Range("C3").Select()

# --- Add Sheet2 after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws2.Columns.Item(1).ColumnWidth = 58.140625

$ws2.Range("A1").Value = "part_name"
$ws2.Range("B1").Value = "no_work _order"

$a2 = $ws2.Range("A2")
$a2.Value = "GARNISH RR BUMPER LWR (IPL)"
$a2.Font.Name = "Arial"
$a2.Font.Size = 10
$a2.Font.Color = 7697781

$b2 = $ws2.Range("B2")
$b2.Value = "23102022SPK001"
$b2.Font.Name = "Arial"
$b2.Font.Size = 10

$ws2.Range("B2").Select()

# restore active sheet to Sheet1
$ws1.Activate()
